$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exported")

# Insert a new column before column A, shifting the existing A/B/C -> B/C/D
$ws.Range("A1:A37").EntireColumn.Insert()

# Header for the new ID column
$ws.Range("A1").Value = "ID"

# Animal IDs, copied from Sheet1 column A (rows 18-53), one per Exported row (2-37)
$ids = @(3035,3036,3037,3038,3039,3040,3041,3042,3048,3049,3050,3055,3056,3057,3058,3047,1747,1749,1750,1766,1763,2102,2098,2099,2105,2107,1748,1751,1745,1746,2108,2111,2118,2123,2127,2128)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
